$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "CRA 97a # 2 - 70"
$ws.Cells.Item(2, 1).Value = "CRA 97a # 2 70"
$ws.Cells.Item(3, 1).Value = "CRA 97a No 2 - 70"
$ws.Cells.Item(4, 1).Value = "CRA 97a No 2 70"
$ws.Cells.Item(5, 1).Value = "CRA 97a Num 2 - 70"
$ws.Cells.Item(6, 1).Value = "CRA 97a Num 2 70"
$ws.Cells.Item(7, 1).Value = "CRA 97a Numero 2 - 70"
$ws.Cells.Item(8, 1).Value = "CRA 97a Numero 2 70"
$ws.Cells.Item(9, 1).Value = "Carrera 97a # 2 - 70"
$ws.Cells.Item(10, 1).Value = "Carrera 97a # 2 70"
$ws.Cells.Item(11, 1).Value = "Carrera 97a No 2 - 70"
$ws.Cells.Item(12, 1).Value = "Carrera 97a No 2 70"
$ws.Cells.Item(13, 1).Value = "Carrera 97a Num 2 - 70"
$ws.Cells.Item(14, 1).Value = "Carrera 97a Num 2 70"
$ws.Cells.Item(15, 1).Value = "Carrera 97a Numero 2 - 70"
$ws.Cells.Item(16, 1).Value = "Carrera 97a Numero 2 70"
$ws.Cells.Item(17, 1).Value = "Calle 97a # 2 - 70"
$ws.Cells.Item(18, 1).Value = "Calle 97a # 2 70"
$ws.Cells.Item(19, 1).Value = "Calle 97a No 2 - 70"
$ws.Cells.Item(20, 1).Value = "Calle 97a No 2 70"
$ws.Cells.Item(21, 1).Value = "Calle 97a Num 2 - 70"
$ws.Cells.Item(22, 1).Value = "Calle 97a Num 2 70"
$ws.Cells.Item(23, 1).Value = "Calle 97a Numero 2 - 70"
$ws.Cells.Item(24, 1).Value = "Calle 97a Numero 2 70"
$ws.Cells.Item(25, 1).Value = "Cl 97a # 2 - 70"
$ws.Cells.Item(26, 1).Value = "Cl 97a # 2 70"
$ws.Cells.Item(27, 1).Value = "Cl 97a No 2 - 70"
$ws.Cells.Item(28, 1).Value = "Cl 97a No 2 70"
$ws.Cells.Item(29, 1).Value = "Cl 97a Num 2 - 70"
$ws.Cells.Item(30, 1).Value = "Cl 97a Num 2 70"
$ws.Cells.Item(31, 1).Value = "Cl 97a Numero 2 - 70"
$ws.Cells.Item(32, 1).Value = "Cl 97a Numero 2 70"
$ws.Cells.Item(33, 1).Value = "Transversal 97a # 2 - 70"
$ws.Cells.Item(34, 1).Value = "Transversal 97a # 2 70"
$ws.Cells.Item(35, 1).Value = "Transversal 97a No 2 - 70"
$ws.Cells.Item(36, 1).Value = "Transversal 97a No 2 70"
$ws.Cells.Item(37, 1).Value = "Transversal 97a Num 2 - 70"
$ws.Cells.Item(38, 1).Value = "Transversal 97a Num 2 70"
$ws.Cells.Item(39, 1).Value = "Transversal 97a Numero 2 - 70"
$ws.Cells.Item(40, 1).Value = "Transversal 97a Numero 2 70"
$ws.Cells.Item(41, 1).Value = "Tv 97a # 2 - 70"
$ws.Cells.Item(42, 1).Value = "Tv 97a # 2 70"
$ws.Cells.Item(43, 1).Value = "Tv 97a No 2 - 70"
$ws.Cells.Item(44, 1).Value = "Tv 97a No 2 70"
$ws.Cells.Item(45, 1).Value = "Tv 97a Num 2 - 70"
$ws.Cells.Item(46, 1).Value = "Tv 97a Num 2 70"
$ws.Cells.Item(47, 1).Value = "Tv 97a Numero 2 - 70"
$ws.Cells.Item(48, 1).Value = "Tv 97a Numero 2 70"

$ws.Range("A1:A48").Select()
